# Daily attendance processing - 2025-12-15 10:33:34
# Reorders the email addresses listed in the "Recorded By" (column G)
# cells for a number of rows on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = "mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    4  = "hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    5  = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    6  = "nourhan.mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    7  = "mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    8  = "backup@backdoor.com, eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    9  = "mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    10 = "drsarahmw@gmail.com, rana.abozaid@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, heba@med.asu.edu.eg"
    11 = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    12 = "System, salma.elgendy.std@med.asu.edu.eg"
    13 = "Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
    15 = "wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
    16 = "norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
    17 = "yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
    19 = "mariam.e.mohammed@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
    20 = "user@user.com, nourhan.mostafa@med.asu.edu.eg"
    22 = "nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
    23 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    24 = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
    25 = "Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    26 = "gigi@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
    27 = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
    28 = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
    29 = "Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
    30 = "hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    31 = "mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    32 = "hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    33 = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hananragab@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    34 = "nourhan.mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    35 = "mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    36 = "backup@backdoor.com, eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    37 = "mennatulla.medhat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    38 = "drsarahmw@gmail.com, rana.abozaid@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, heba@med.asu.edu.eg"
    39 = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    40 = "System, salma.elgendy.std@med.asu.edu.eg"
    41 = "Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
    43 = "wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
    44 = "norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
    45 = "yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
    47 = "mariam.e.mohammed@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
    48 = "user@user.com, nourhan.mostafa@med.asu.edu.eg"
    50 = "nancy.abdelshafy@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
    51 = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
    52 = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
    53 = "Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    54 = "gigi@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
    55 = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
    56 = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
    57 = "Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
}

foreach ($rowNum in $newValues.Keys) {
    $ws.Cells.Item($rowNum, 7).Value = $newValues[$rowNum]
}
